# Append the latest tracker snapshot (2025-09-11) to the progress history sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 45911
$progress = 0.9609803444828162

$goals = @(
    @{ Id = "G2"; Name = "Workout" },
    @{ Id = "G3"; Name = "Eat Healthy" },
    @{ Id = "G4"; Name = "Read Book" },
    @{ Id = "G5"; Name = "Investment Plan" },
    @{ Id = "G6"; Name = "Spend 10 Hours without phone" }
)

$startRow = 22
for ($i = 0; $i -lt $goals.Count; $i++) {
    $row = $startRow + $i
    $goal = $goals[$i]

    $ws.Cells.Item($row, 1).Value = $goal.Id
    $ws.Cells.Item($row, 2).Value = $goal.Name
    $ws.Cells.Item($row, 3).Value = $newDate
    $ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item(2, 3).NumberFormat
    $ws.Cells.Item($row, 4).Value = $progress
    $ws.Cells.Item($row, 5).Value = 0
    $ws.Cells.Item($row, 6).Value = -0.01
}
